$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the data held in rows 4, 5 and 6:
#   new row 4 <- old row 6 (Tretåig hackspett / Picoides tridactylus)
#   new row 5 <- old row 4 (Gammelgransskål / Pseudographis pinicola)
#   new row 6 <- old row 5 (Garnlav / Alectoria sarmentosa)
# Columns that hold identical values across all three rows (D, I, K, N, P, S,
# T, U, V, W, Y, AA, AD, AG, AT, AW, AX, AY) are left untouched.

# ---- Row 4 (becomes the former row 6 record) ----
$ws.Range("A4").Value = 131136961
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("J4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("Q4").Value = 789068
$ws.Range("R4").Value = 7131245
$ws.Range("AC4").Value = "barksprätt på gammal gran"
$ws.Range("AE4").Value = $true
$ws.Range("AF4").Value = ""

# ---- Row 5 (becomes the former row 4 record) ----
$ws.Range("A5").Value = 131136941
$ws.Range("B5").Value = 83089
$ws.Range("E5").Value = 1312
$ws.Range("F5").Value = "Gammelgransskål"
$ws.Range("G5").Value = "Pseudographis pinicola"
$ws.Range("H5").Value = "(Nyl.) Rehm"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("Q5").Value = 788995
$ws.Range("R5").Value = 7131220
$ws.Range("AC5").Value = "på en gammal senvuxen gran"
$ws.Range("AE5").Value = $false

# ---- Row 6 (becomes the former row 5 record) ----
$ws.Range("A6").Value = 131136874
$ws.Range("B6").Value = 79243
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("J6").Value = "bålar"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 788960
$ws.Range("R6").Value = 7131416
$ws.Range("AC6").Value = ""
$ws.Range("AE6").Value = $false
$ws.Range("AF6").Value = ""
